$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 1331.7715
$ws.Range("I137").Value = 1067.6538
$ws.Range("J137").Value = 2094.7778
$ws.Range("K137").Value = 3202.9614
$ws.Range("L137").Value = 6284.3334
$ws.Range("M137").Value = -652.9614000000001
$ws.Range("N137").Value = -11384.3334

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 3263.568
$ws.Range("I61").Value = 3323.2563
$ws.Range("J61").Value = 2798
$ws.Range("K61").Value = 3323.2563
$ws.Range("L61").Value = 2798
$ws.Range("M61").Value = -3111.2563
$ws.Range("N61").Value = -3222
$ws.Range("H74").Value = 5452.8335
$ws.Range("I74").Value = 991.75
$ws.Range("K74").Value = 991.75
$ws.Range("M74").Value = -117.75
$ws.Range("H77").Value = 5452.8335
$ws.Range("I77").Value = 991.75
$ws.Range("K77").Value = 4958.75
$ws.Range("M77").Value = -590.75
$ws.Range("H102").Value = 111112800
$ws.Range("I102").Value = 125001650
$ws.Range("J102").Value = 2000
$ws.Range("K102").Value = 125001650
$ws.Range("L102").Value = 2000
$ws.Range("M102").Value = -125000028
$ws.Range("N102").Value = -5244
$ws.Range("H132").Value = 3789362.5
$ws.Range("I132").Value = 5682802.5
$ws.Range("J132").Value = 2482.7273
$ws.Range("K132").Value = 17048407.5
$ws.Range("L132").Value = 7448.1819
$ws.Range("M132").Value = -17045877.5
$ws.Range("N132").Value = -12508.1819
$ws.Range("H136").Value = 3263.568
$ws.Range("I136").Value = 3323.2563
$ws.Range("J136").Value = 2798
$ws.Range("K136").Value = 9969.768899999999
$ws.Range("L136").Value = 8394
$ws.Range("M136").Value = -7419.768899999999
$ws.Range("N136").Value = -13494

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 1153.6
$ws.Range("J86").Value = 960
$ws.Range("L86").Value = 960
$ws.Range("N86").Value = -3206
$ws.Range("H89").Value = 1153.6
$ws.Range("J89").Value = 960
$ws.Range("L89").Value = 4800
$ws.Range("N89").Value = -16032
$ws.Range("H99").Value = 1072.8572
$ws.Range("I99").Value = 1001.6667
$ws.Range("K99").Value = 1001.6667
$ws.Range("M99").Value = 496.3333
$ws.Range("H103").Value = 22689.834
$ws.Range("J103").Value = 22689.834
$ws.Range("L103").Value = 22689.834
$ws.Range("N103").Value = -25033.834
$ws.Range("H105").Value = 1604.7142
$ws.Range("J105").Value = 1622.1666
$ws.Range("L105").Value = 1622.1666
$ws.Range("N105").Value = -5116.1666
$ws.Range("H107").Value = 3367.625
$ws.Range("I107").Value = 3785.25
$ws.Range("K107").Value = 3785.25
$ws.Range("M107").Value = -1865.25
$ws.Range("H134").Value = 20744.654
$ws.Range("I134").Value = 21374.44
$ws.Range("J134").Value = 5000
$ws.Range("K134").Value = 64123.31999999999
$ws.Range("L134").Value = 15000
$ws.Range("M134").Value = -61588.31999999999
$ws.Range("N134").Value = -20070

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 8818.642
$ws.Range("I31").Value = 926.24243
$ws.Range("J31").Value = 21841.1
$ws.Range("K31").Value = 926.24243
$ws.Range("L31").Value = 21841.1
$ws.Range("M31").Value = -631.24243
$ws.Range("N31").Value = -22431.1
$ws.Range("H34").Value = 8818.642
$ws.Range("I34").Value = 926.24243
$ws.Range("J34").Value = 21841.1
$ws.Range("K34").Value = 926.24243
$ws.Range("L34").Value = 21841.1
$ws.Range("M34").Value = -724.24243
$ws.Range("N34").Value = -22245.1
$ws.Range("H58").Value = 3694706
$ws.Range("I58").Value = 5755106
$ws.Range("K58").Value = 5755106
$ws.Range("M58").Value = -5754903
$ws.Range("H132").Value = 16673752
$ws.Range("I132").Value = 47620160
$ws.Range("J132").Value = 10300.923
$ws.Range("K132").Value = 142860480
$ws.Range("L132").Value = 30902.769
$ws.Range("M132").Value = -142857950
$ws.Range("N132").Value = -35962.769
$ws.Range("H134").Value = 8447189
$ws.Range("I134").Value = 9616482
$ws.Range("J134").Value = 5683407
$ws.Range("K134").Value = 28849446
$ws.Range("L134").Value = 17050221
$ws.Range("M134").Value = -28846911
$ws.Range("N134").Value = -17055291
$ws.Range("H136").Value = 3694706
$ws.Range("I136").Value = 5755106
$ws.Range("K136").Value = 17265318
$ws.Range("M136").Value = -17262768

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H140").Value = 1091.9269
$ws.Range("I140").Value = 737.3
$ws.Range("J140").Value = 2059.0908
$ws.Range("K140").Value = 2211.9
$ws.Range("L140").Value = 6177.2724
$ws.Range("M140").Value = 2968.1
$ws.Range("N140").Value = -16537.2724

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 50005790
$ws.Range("I132").Value = 125001064
$ws.Range("J132").Value = 8943.833000000001
$ws.Range("K132").Value = 375003192
$ws.Range("L132").Value = 26831.499
$ws.Range("M132").Value = -375000662
$ws.Range("N132").Value = -31891.499

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H93").Value = 85048.414
$ws.Range("I93").Value = 554.4286
$ws.Range("J93").Value = 203340
$ws.Range("K93").Value = 554.4286
$ws.Range("L93").Value = 203340
$ws.Range("M93").Value = 693.5714
$ws.Range("N93").Value = -205836
$ws.Range("H132").Value = 11113770
$ws.Range("I132").Value = 20002188
$ws.Range("J132").Value = 3245.5
$ws.Range("K132").Value = 60006564
$ws.Range("L132").Value = 9736.5
$ws.Range("M132").Value = -60004034
$ws.Range("N132").Value = -14796.5
$ws.Range("H136").Value = 4485.725
$ws.Range("I136").Value = 5414.448
$ws.Range("K136").Value = 16243.344
$ws.Range("M136").Value = -13693.344

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H2").Value = 167590800
$ws.Range("I2").Value = 1800
$ws.Range("J2").Value = 201108600
$ws.Range("K2").Value = 1800
$ws.Range("L2").Value = 201108600
$ws.Range("M2").Value = -1688
$ws.Range("N2").Value = -201108824
$ws.Range("H132").Value = 63486860
$ws.Range("I132").Value = 57147076
$ws.Range("J132").Value = 85676110
$ws.Range("K132").Value = 171441228
$ws.Range("L132").Value = 257028330
$ws.Range("M132").Value = -171438698
$ws.Range("H136").Value = 18542394
$ws.Range("I136").Value = 14538008
$ws.Range("J136").Value = 22728798
$ws.Range("K136").Value = 43614024
$ws.Range("L136").Value = 68186394
$ws.Range("M136").Value = -43611474
